$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Change all cells that previously held "CHO_DUYET" to "DA_DUYET"
$ws.Range("H2").Value = "DA_DUYET"
$ws.Range("H3").Value = "DA_DUYET"

# H4 is now a separate, new status: "HUY"
$ws.Range("H4").Value = "HUY"
